$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Nome do Plano" -> "Nome do " + "Cliente" (two runs)
# ---------------------------------------------------------------
$countP = $d.Paragraphs.Count
$idxNome = -1
for ($i = 1; $i -le $countP; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Nome do Plano*") {
        $idxNome = $i
        break
    }
}
$pNome = $d.Paragraphs.Item($idxNome)
$pNomeStart = $pNome.Range.Start

# Replace just the word "Plano" with "Cliente" (still a single run after this)
$rPlano = $d.Range($pNomeStart, $pNomeStart + "Nome do Plano".Length)
$rPlano.Find.Execute("Plano", $true, $false, $false, $false, $false, $true, 1, $false, "Cliente", 2)

# Force a run split exactly between "Nome do " and "Cliente" by toggling a
# character-formatting property on just the "Cliente" portion and back off;
# this makes the engine materialize a distinct <w:r> (with unchanged rPr)
# instead of merging it back with its neighbour.
$splitStart = $pNomeStart + "Nome do ".Length
$splitEnd = $pNomeStart + "Nome do Cliente".Length
$rSplit1 = $d.Range($splitStart, $splitEnd)
$rSplit1.Bold = 1
$rSplit1.Bold = 0

# ---------------------------------------------------------------
# 2) "Tipo do Plano" -> "Nível" + " do Plano" (two runs)
# ---------------------------------------------------------------
$idxTipo = -1
for ($i = 1; $i -le $countP; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Tipo do Plano*") {
        $idxTipo = $i
        break
    }
}
$pTipo = $d.Paragraphs.Item($idxTipo)
$pTipoStart = $pTipo.Range.Start

$rTipo = $d.Range($pTipoStart, $pTipoStart + "Tipo do Plano".Length)
$rTipo.Find.Execute("Tipo", $true, $false, $false, $false, $false, $true, 1, $false, "Nível", 2)

$splitStart2 = $pTipoStart + "Nível".Length
$splitEnd2 = $pTipoStart + "Nível do Plano".Length
$rSplit2 = $d.Range($splitStart2, $splitEnd2)
$rSplit2.Bold = 1
$rSplit2.Bold = 0

# ---------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the "Data" paragraph
#    to the end of the (now) "Nível do Plano" paragraph.
# ---------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# Re-resolve the "Nível do Plano" paragraph end (text length changed above).
$pTipo2 = $d.Paragraphs.Item($idxTipo)
$paraEnd = $pTipo2.Range.End - 1

# Adding a bookmark exactly at a paragraph-end offset is mishandled by this
# host, so nudge the boundary: insert a throwaway character after the
# target offset, anchor the bookmark just before it, then remove the
# throwaway character again (the bookmark collapses back correctly).
$guard = $d.Range($paraEnd, $paraEnd)
$guard.InsertAfter("#")

$bmRange = $d.Range($paraEnd, $paraEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$guard2 = $d.Range($paraEnd, $paraEnd + 1)
$guard2.Delete()
